# Fix a typo in the "keyword" column of the "test steps" sheet:
# the shared string "fillTextFiled" should read "fillTextField".
# Update every cell that uses it (F4, F5, F10), then leave the
# selection on the last-touched cell (F10), matching the author's
# saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("test steps")

$ws.Range("F4").Value = "fillTextField"
$ws.Range("F5").Value = "fillTextField"
$ws.Range("F10").Value = "fillTextField"

$ws.Activate() | Out-Null
$ws.Range("F10").Select() | Out-Null
